$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header row (date label, Ballgorithm, ESPN)
$ws.Range("A1").Value = " NBA, Sunday 11th Feb 2024"
$ws.Range("B1").Value = "Ballgorithm"
$ws.Range("C1").Value = "ESPN"

# Row 2
$ws.Range("A2").Value = "Boston Celtics (40-12) vs Miami Heat (28-24)"
$ws.Range("B2").Value = " Boston Celtics (89.29%)"
$ws.Range("C2").Value = " Boston Celtics (69.5%)"

# Row 3
$ws.Range("A3").Value = "Sacramento Kings (30-21) vs Oklahoma City Thunder (35-17)"
$ws.Range("B3").Value = "Oklahoma City Thunder (76.92%)"
$ws.Range("C3").Value = "Oklahoma City Thunder (64.0%)"

# Remove the old rows 4-12 so the used range shrinks back to A1:C3
$ws.Range("A4:C12").Clear()

# Match the saved selection/active cell from the authored workbook
$null = $ws.Range("A3").Select()
